# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 197
$ws1.Range("F3").Value = 5445
$ws1.Range("F8").Value = 603
$ws1.Range("F12").Value = 4714
$ws1.Range("F17").Value = 3599
$ws1.Range("F19").Value = 1125
$ws1.Range("F20").Value = 111
$ws1.Range("F27").Value = 76
$ws1.Range("F28").Value = 328
$ws1.Range("F29").Value = 37
$ws1.Range("F30").Value = 62
$ws1.Range("F32").Value = 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 197
$ws4.Range("F4").Value = 5445
$ws4.Range("F9").Value = 603
$ws4.Range("F13").Value = 4714
$ws4.Range("F18").Value = 3599
$ws4.Range("F20").Value = 1125
$ws4.Range("F21").Value = 111
$ws4.Range("F28").Value = 76
$ws4.Range("F29").Value = 328
$ws4.Range("F30").Value = 37
$ws4.Range("F31").Value = 62
$ws4.Range("F33").Value = 34
